# Apply updated cryptocurrency price/volume figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values are plain text (e.g. "58.645.07", "23.70",
# "0.0₃09759") and must stay text, so each is written with a leading
# apostrophe -- Excel's standard "force text" marker -- to stop values
# like "23.70" from being silently coerced into the number 23.7 (which
# would drop the trailing zero) or otherwise reformatted as a number/date.
# Column E ("Volume(1h)") values already carry surrounding whitespace
# (e.g. "  -6.69%  "), so Excel keeps those as text with no marker needed.

$ws.Range("D2").Value = '''58.645.07'
$ws.Range("E2").Value = '  -6.69%  '
$ws.Range("D3").Value = '''2.432.26'
$ws.Range("E3").Value = '  -9.33%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''533.07'
$ws.Range("E5").Value = '  -4.04%  '
$ws.Range("D6").Value = '''144.62'
$ws.Range("E6").Value = '  -8.31%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '''0.571'
$ws.Range("E8").Value = '  -3.11%  '
$ws.Range("D9").Value = '''2.445.15'
$ws.Range("E9").Value = '  -8.98%  '
$ws.Range("D10").Value = '''0.0984'
$ws.Range("E10").Value = '  -7.19%  '
$ws.Range("E11").Value = '  -2.23%  '
$ws.Range("D12").Value = '''5.28'
$ws.Range("E12").Value = '  -2.15%  '
$ws.Range("D13").Value = '''0.348'
$ws.Range("E13").Value = '  -5.38%  '
$ws.Range("D14").Value = '''2.872.81'
$ws.Range("E14").Value = '  -9.00%  '
$ws.Range("D15").Value = '''23.70'
$ws.Range("E15").Value = '  -10.54%  '
$ws.Range("D16").Value = '''58.552.54'
$ws.Range("E16").Value = '  -6.69%  '
$ws.Range("D17").Value = '''0.0000137'
$ws.Range("E17").Value = '  -6.32%  '
$ws.Range("D18").Value = '''2.490.30'
$ws.Range("E18").Value = '  -7.24%  '
$ws.Range("D19").Value = '''11.03'
$ws.Range("E19").Value = '  -6.68%  '
$ws.Range("D20").Value = '''4.32'
$ws.Range("E20").Value = '  -6.39%  '
$ws.Range("D21").Value = '''321.68'
$ws.Range("E21").Value = '  -6.78%  '
$ws.Range("E22").Value = '  -3.51%  '
$ws.Range("D23").Value = '''5.65'
$ws.Range("E23").Value = '  -9.11%  '
$ws.Range("D24").Value = '''60.34'
$ws.Range("E24").Value = '  -4.59%  '
$ws.Range("D25").Value = '''0.446'
$ws.Range("E25").Value = '  -12.67%  '
$ws.Range("E26").Value = '  -6.41%  '
$ws.Range("D27").Value = '''0.975'
$ws.Range("E27").Value = '  -2.42%  '
$ws.Range("D28").Value = '''7.59'
$ws.Range("E28").Value = '  -7.08%  '
$ws.Range("D29").Value = '''1.80'
$ws.Range("E29").Value = '  -7.25%  '
$ws.Range("D30").Value = '''0.0₃0759'
$ws.Range("E30").Value = '  -10.98%  '
$ws.Range("D31").Value = '''6.58'
$ws.Range("E31").Value = '  -9.20%  '
$ws.Range("D32").Value = '''1.19'
$ws.Range("E32").Value = '  -13.87%  '
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("D34").Value = '''155.55'
$ws.Range("E34").Value = '  -5.46%  '
$ws.Range("D35").Value = '''18.40'
$ws.Range("E35").Value = '  -5.56%  '
$ws.Range("D36").Value = '''1.36'
$ws.Range("E36").Value = '  -7.71%  '
$ws.Range("D37").Value = '''4.37'
$ws.Range("E37").Value = '  -10.87%  '
$ws.Range("D38").Value = '''1.65'
$ws.Range("E38").Value = '  -6.99%  '
$ws.Range("D39").Value = '''5.71'
$ws.Range("E39").Value = '  -7.65%  '
$ws.Range("D40").Value = '''307.64'
$ws.Range("E40").Value = '  -11.03%  '
$ws.Range("D41").Value = '''36.12'
$ws.Range("E41").Value = '  -5.82%  '
$ws.Range("D42").Value = '''0.822'
$ws.Range("E42").Value = '  -12.23%  '
$ws.Range("D43").Value = '''3.67'
$ws.Range("E43").Value = '  -7.95%  '
$ws.Range("D44").Value = '''0.999'
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("E45").Value = '  -2.29%  '
$ws.Range("D46").Value = '''0.0935'
$ws.Range("E46").Value = '  -3.58%  '
$ws.Range("D47").Value = '''0.577'
$ws.Range("E47").Value = '  -6.63%  '
$ws.Range("D48").Value = '''0.0521'
$ws.Range("E48").Value = '  -6.27%  '
$ws.Range("D49").Value = '''0.0227'
$ws.Range("E49").Value = '  -5.89%  '
$ws.Range("D50").Value = '''121.02'
$ws.Range("E50").Value = '  -6.03%  '
$ws.Range("D51").Value = '''18.32'
$ws.Range("E51").Value = '  -9.11%  '
